# Add entrance animation ("Appear") to the system/architecture diagram
# slide: each "Flowchart: Document" node appears on click, and its
# connecting arrow appears together with it ("with previous").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

function Get-ShapeById($slide, $id) {
    foreach ($sh in $slide.Shapes) {
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# msoAnimEffectAppear = 1
# msoAnimTriggerWithPrevious = 2
function Add-AppearPair($clickShapeId, $withShapeId) {
    $clickShape = Get-ShapeById $s $clickShapeId
    $withShape = Get-ShapeById $s $withShapeId

    $null = $s.TimeLine.MainSequence.AddEffect($clickShape, 1)
    $null = $s.TimeLine.MainSequence.AddEffect($withShape, 1, 0, 2)
}

Add-AppearPair 24 21
Add-AppearPair 25 10
Add-AppearPair 26 11
Add-AppearPair 27 37
Add-AppearPair 36 33
Add-AppearPair 29 28
Add-AppearPair 38 40
